# Update "想去人数" (want-to-go count) figures across sheets to reflect
# the latest scrape of gh-pages output (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition) ---
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F5").Value  = 749
$wsExhibit.Range("F6").Value  = 2354
$wsExhibit.Range("F9").Value  = 2965
$wsExhibit.Range("F11").Value = 4424
$wsExhibit.Range("F12").Value = 387
$wsExhibit.Range("F13").Value = 212
$wsExhibit.Range("F18").Value = 209
$wsExhibit.Range("F22").Value = 4471
$wsExhibit.Range("F25").Value = 1138
$wsExhibit.Range("F27").Value = 564
$wsExhibit.Range("F30").Value = 579
$wsExhibit.Range("F32").Value = 539

# --- Sheet "本地生活" (Local Life) ---
$wsLocal = $wb.Worksheets.Item("本地生活")
$wsLocal.Range("F3").Value = 1035

# --- Sheet "全部类型" (All Types) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value  = 1035
$wsAll.Range("F8").Value  = 749
$wsAll.Range("F9").Value  = 2354
$wsAll.Range("F13").Value = 2965
$wsAll.Range("F15").Value = 4424
$wsAll.Range("F16").Value = 387
$wsAll.Range("F17").Value = 212
$wsAll.Range("F22").Value = 209
$wsAll.Range("F27").Value = 4471
$wsAll.Range("F30").Value = 1138
$wsAll.Range("F32").Value = 564
$wsAll.Range("F35").Value = 579
$wsAll.Range("F37").Value = 539
